$d = $word.ActiveDocument

$replacements = @(
    @{old="688×8=5504"; new="300×3=900"},
    @{old="306×5=1530"; new="710×5=3550"},
    @{old="245×2=490"; new="423×2=846"},
    @{old="931×3=2793"; new="108×9=972"},
    @{old="847×5=4235"; new="971×5=4855"},
    @{old="135×3=405"; new="238×6=1428"},
    @{old="945×2=1890"; new="114×6=684"},
    @{old="709×3=2127"; new="973×3=2919"},
    @{old="246×9=2214"; new="192×4=768"},
    @{old="503×4=2012"; new="241×3=723"},
    @{old="487×7=3409"; new="769×7=5383"},
    @{old="375×5=1875"; new="348×8=2784"},
    @{old="635×6=3810"; new="415×7=2905"},
    @{old="593×6=3558"; new="688×7=4816"},
    @{old="105×8=840"; new="793×5=3965"},
    @{old="191×9=1719"; new="978×9=8802"},
    @{old="848×9=7632"; new="458×7=3206"},
    @{old="521×9=4689"; new="456×8=3648"},
    @{old="101×3=303"; new="585×6=3510"},
    @{old="240×7=1680"; new="582×4=2328"},
    @{old="103×7=721"; new="766×4=3064"},
    @{old="930×2=1860"; new="908×3=2724"},
    @{old="545×4=2180"; new="147×4=588"},
    @{old="930×3=2790"; new="552×6=3312"},
    @{old="563×6=3378"; new="654×3=1962"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}

Write-Output "Done"
